$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new (blank) column at M ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column M; existing M (In Advance), N (Late),
# O (heading) and P (Outstanding) shift right to N, O, P, Q respectively.
$ws.Columns("M:M").Insert() | Out-Null

# Give the newly inserted column the same width as column K (the post-edit
# workbook shows the new column with that width as a manual/custom width,
# i.e. without the "best fit" flag the other columns carry).
$ws.Columns("M:M").ColumnWidth = $ws.Columns("K:K").ColumnWidth

# Update the repayment schedule numbers: the amount that was "Due" moved
# into the new "In Advance" position while Due itself became 0.
$ws.Range("K3").Value = 0

$ws.Range("N3").ClearFormats()
$ws.Range("N3").NumberFormat = "#,##0"
$ws.Range("N3").Value = 10000

# Make "Repayment schedule" the active sheet / tab, with K9 selected.
# (Do this last: activating/selecting on another sheet would otherwise
# steal the "active sheet" status back.)
$ws.Activate() | Out-Null
$ws.Range("K9").Select() | Out-Null
